$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old `_GoBack` bookmark that currently sits after
#    "What are the techniques?" - it moves elsewhere in this edit.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "What is the weekly report?" -> "What is the weekly report
#    (written or vocal)?", split across four runs with identical
#    formatting (b, sz 28).
# ---------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("What is the weekly report?", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$start = $find.Start

$find.Text = "What is the weekly report (written or vocal)?"

# Split the single run produced by the text assignment into four runs
# by briefly adding (and immediately removing) a bookmark at each
# desired boundary - the bookmark forces a run break that survives its
# own deletion.
$splitPoints = @(19, 27, 44)
$i = 0
foreach ($offset in $splitPoints) {
    $i = $i + 1
    $pos = $start + $offset
    $bmName = "TempSplitA$i"
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($bmName, $r)
    $d.Bookmarks($bmName).Delete()
}

# ---------------------------------------------------------------------
# 3) "How comprehensive does the sponsor want this chapter?" - split
#    "want" into "wan" / "t" and drop the `_GoBack` bookmark in between.
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("How comprehensive does the sponsor wan", $true, $false, $false, $false, `
                     $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $find2.End
$r2 = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $r2)
